$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parts - Consoles.csv")
$ws2 = $wb.Worksheets.Item("revision")

$ws.Range('AG1').Value = 'Incline Calibration'
$ws.Range('AH1').Value = 'Display Test'
$ws.Range('AI1').Value = 'Display Test Image'
$ws.Range('AJ1').Value = 'Button Test'
$ws.Range('AK1').Value = 'Button Test Image'
$ws.Range('AL1').Value = 'Drive Motor Output Test'
$ws.Range('AM1').Value = 'Tach Input Test'
$ws.Range('AN1').Value = 'Resistance Motor Test'
$ws.Range('AO1').Value = 'Incline Motor Test'
$ws.Range('AP1').Value = 'USB port test'
$ws.Range('AQ1').Value = 'Hand Pulse Test'
$ws.Range('AR1').Value = 'Chest Pulse Test'
$ws.Range('AS1').Value = 'Fan Test'
$ws.Range('AT1').Value = 'Audio Test'
$ws.Range('AU1').Value = 'TV Test'
$ws.Range('AV1').Value = 'Upright Motor Test'
$ws.Range('AW1').Value = 'Finish Test'
$ws.Range('AX1').Value = 'Last User'
$ws.Range('AY1').Value = 'Created By User'
$ws.Range('AZ1').Value = 'Updated'
$ws.Range('AK2').ClearContents()
$ws.Range('AN2').ClearContents()
$ws.Range('AO2').Value = 'On the console, press random Quick incline buttons.
Verify that the Incline on the console and EQF1259 match.'
$ws.Range('AP2').Value = 'Plug USB cable from EQF1259 into console and verify "USB" is displayed on the EQF1259.'
$ws.Range('AQ2').Value = 'Hold the pulse bars and verify a pulse reading is displayed'
$ws.Range('AR2').Value = 'Press Display button to goto BLE Pulse screen.
Press Start button.
Run manual workout on the console.
Verify pulse BLE pulse is read on console.'
$ws.Range('AS2').Value = 'Connect the 3-pin fans to the console.
On the console, press the Large Fan button to turn fan on. Fan should run on low.
Press the Large Fan button again and the fan should run on high.
Press the Small Fan button to turn fan to low, then press again to turn fan off.'
$ws.Range('AT2').Value = 'Connect the EQW1007 to the iOS audio source and connect the other end of the EQW1007 to the console.
Play audio and verify that it plays out of the console speakers.
Adjust the volume to minimum and maximum level and verify a change.'
$ws.Range('AX2').Value = 'PIP_GEN_ID-0'
$ws.Range('AE2').Value = 'Touch user profile icon at the lower right  of the screen
Touch settings
Touch “Equipment Info” and then “App Info”. This will display App version and Brainboard version.
Verify and Record.
Press the back arrow on the tablet to return to main screen'
$ws.Range('AE2').WrapText = $true
$ws.Range('AG2').Value = 'Touch the user profile icon on lower right of the screen
Touch Settings
Touch Maintainance
Select Calibrate Incline
Press Calibrate on EQF1259
Press Begin on tablet to start incline calibration.  Verify the incline value on the display increases briefly, pauses, then decrease to zero
When calibration is complete, press the back arrow on the tablet to return to main screen'
$ws.Range('AG2').WrapText = $true
$ws.Range('AH2').Value = 'Verify tablet display is lit and without flaws'
$ws.Range('AH2').WrapText = $true
$ws.Range('AL2').Value = 'On the console, press random Quick speed buttons.
Verify that the Speed on the console and EQF1259 match.'
$ws.Range('AL2').WrapText = $true
$ws.Range('AU2').Value = 'Connect the HDMI cable to the TV.
Connect the HDMI cable to the console and the console screen and sound will be mirrored on the TV.'
$ws.Range('AU2').WrapText = $true
$ws.Range('AZ2').Value = 43119.5863994213
$ws.Range('AZ2').NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws2.Range('D1').Value = 'id-0'
$ws2.Range('C2').Value = 'Created arrays for variables'
$ws2.Range('D2').Value = 'id-0'
$ws2.Range('C3').Value = 'Added time stamping'
$ws2.Range('D3').Value = 'id-0'
$ws2.Range('D4').Value = 'id-0'
$ws2.Range('C5').Value = 'Converted to openpyxl'
$ws2.Range('D5').Value = 'id-0'
$ws2.Range('A6').Value = 'REV 5'
$ws2.Range('C6').Value = 'Killed trackvia full table'
$ws2.Range('D6').Value = 'id-0'
$ws2.Range('B7').Value = 43119.5863697338
$ws2.Range('D7').Value = 'id-0'

# revision sheet column C width change (79.2 -> 56.4)
$ws2.Columns.Item(3).ColumnWidth = 56.4
